$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers in row 1 to clean, machine-friendly column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization of connector words ("de"/"del"/"los" -> "De"/"Del"/"Los")
# in municipality / state names
$ws.Range("B3").Value = "Amatenango De La Frontera"
$ws.Range("A11").Value = "Ciudad De México"
$ws.Range("A20").Value = "Estado De México"
$ws.Range("B24").Value = "Tlalnepantla De Baz"
$ws.Range("B27").Value = "Valle De Chalco Solidaridad"
$ws.Range("B30").Value = "Acapulco De Juárez"
$ws.Range("B33").Value = "Ayutla De Los Libres"
$ws.Range("B34").Value = "Chilapa De Álvarez"
$ws.Range("B35").Value = "Chilpancingo De Los Bravo"
$ws.Range("B38").Value = "Coyuca De Benítez"
$ws.Range("B41").Value = "Tlapa De Comonfort"
$ws.Range("B53").Value = "Zacualpan De Amilpas"
$ws.Range("B62").Value = "San Miguel Del Puerto"
$ws.Range("B65").Value = "Villa De Zaachila"
$ws.Range("B68").Value = "Chalchicomula De Sesma"
$ws.Range("B71").Value = "Izúcar De Matamoros"
$ws.Range("B78").Value = "Tepanco De López"
$ws.Range("B80").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B81").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B98").Value = "Cosamaloapan De Carpio"
$ws.Range("B100").Value = "Paso Del Macho"

# Remove trailing footer/metadata rows (110-114), keeping row 109 blank as-is
$ws.Range("A110:D114").EntireRow.Delete()
